{"js": "// Word JS API (Office.js) edit script.\n// Body of: async (context) => { ... }\n//\n// Target change (from the OOXML diff):\n//  1. Remove the old \"_GoBack\" bookmark that sits by itself in an empty\n//     paragraph further down the document (right after the paragraph\n//     ending in \"...oppl\u00e6ringen.\").\n//  2. Delete the whole paragraph that contains the placeholder text\n//     \"{utdanningsProgram}\" (it directly follows the paragraph that\n//     contains \"{klasseTrinn}\").\n//  3. Re-create the \"_GoBack\" bookmark at the very end of the\n//     \"{klasseTrinn}\" paragraph (the position the deleted paragraph used\n//     to start at) so a \"Go back to last edit\" still resolves sensibly.\n//\n// Bookmark ids in the saved XML are reassigned sequentially by the\n// document writer, so we don't need to manage numeric ids ourselves \u2014\n// only bookmark names and relative order/position matter.\n\n// Step 1: drop the stray \"_GoBack\" bookmark before we add a new one, so\n// there is never more than one bookmark named \"_GoBack\" at a time (Word\n// itself only ever keeps a single \"_GoBack\" bookmark in a document).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Step 2: locate the two paragraphs involved using the body's paragraph\n// collection (loading \"text\" lets us find them by their placeholder\n// content instead of hard-coded indexes).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet klasseTrinnParagraph = null;\nlet utdanningsProgramParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"{klasseTrinn}\") !== -1) {\n    klasseTrinnParagraph = paragraphs.items[i];\n  } else if (text.indexOf(\"{utdanningsProgram}\") !== -1) {\n    utdanningsProgramParagraph = paragraphs.items[i];\n  }\n}\n\nif (klasseTrinnParagraph && utdanningsProgramParagraph) {\n  // Step 3: insert the \"_GoBack\" bookmark at the end of the\n  // \"{klasseTrinn}\" paragraph (collapsed range, matching what Word does\n  // when it leaves a bookmark mark at the place of the last edit).\n  const endRange = klasseTrinnParagraph.getRange(\"End\");\n  endRange.insertBookmark(\"_GoBack\");\n\n  // Step 4: remove the now-obsolete \"{utdanningsProgram}\" paragraph\n  // entirely.\n  utdanningsProgramParagraph.delete();\n\n  await context.sync();\n}\n", "ps1": "# Word COM interop edit script.\n# $word / $app == Word.Application, $d == $word.ActiveDocument\n#\n# Target change (from the OOXML diff):\n#  1. Remove the old \"_GoBack\" bookmark that sits by itself in an empty\n#     paragraph further down the document (right after the paragraph\n#     ending in \"...oppl\u00e6ringen.\").\n#  2. Delete the whole paragraph that contains the placeholder text\n#     \"{utdanningsProgram}\" (it directly follows the paragraph that\n#     contains \"{klasseTrinn}\").\n#  3. Re-create the \"_GoBack\" bookmark at the very end of the\n#     \"{klasseTrinn}\" paragraph (the position the deleted paragraph used\n#     to start at).\n#\n# Bookmark ids in the saved XML are reassigned sequentially by the\n# document writer, so only bookmark names / relative positions matter,\n# not the numeric w:id values.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: drop the stray \"_GoBack\" bookmark -----------------------\n# Do this before inserting the new one so there is never more than one\n# bookmark named \"_GoBack\" at the same time (real Word only ever keeps a\n# single \"_GoBack\" bookmark in a document).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- Step 2: find the two paragraphs involved -------------------------\n$targetPara = $null\n$removePara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $t = $p.Range.Text\n  if ($t -like \"*{klasseTrinn}*\") {\n    $targetPara = $p\n  } elseif ($t -like \"*{utdanningsProgram}*\") {\n    $removePara = $p\n  }\n}\n\nif ($targetPara -ne $null -and $removePara -ne $null) {\n  # --- Step 3: insert the \"_GoBack\" bookmark right after the visible\n  # text of the \"{klasseTrinn}\" paragraph (i.e. right before its\n  # paragraph mark).\n  #\n  # A collapsed range placed exactly on a paragraph-mark boundary can\n  # get mis-anchored, so we briefly insert a one-character marker right\n  # after the text, collapse the (non-boundary) range in front of the\n  # marker, add the bookmark there, then remove the marker again. The\n  # bookmark itself stays put once created.\n  $anchor = $d.Range($targetPara.Range.End - 1, $targetPara.Range.End - 1)\n  $anchor.InsertAfter([char]1)\n  $anchor.Collapse(1) # wdCollapseStart: in front of the marker char\n  $anchor.Bookmarks.Add(\"_GoBack\")\n\n  $markerRange = $anchor.Duplicate\n  $markerRange.MoveEnd(1, 1)\n  $markerRange.Delete()\n\n  # --- Step 4: remove the now-obsolete \"{utdanningsProgram}\" paragraph\n  # entirely (text + its paragraph mark).\n  $removePara.Range.Delete()\n}\n"}
